# Minor corrections to slide 16 ("Just-In-Time Compiler")
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

# Locate the body placeholder shape ("Rectangle 3") that holds the bullet text.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Rectangle 3") {
        $shp = $candidate
    }
}

$tr = $shp.TextFrame.TextRange

# --- Change 1 -----------------------------------------------------------
# " is a compiler that converts program source code into native machine
# code as the program is running." ->
# " is a compiler that converts program code into native machine code as
# the program is running."
# (drop the word "source " so the run reads "program code" instead of
# "program source code")
$oldRun1 = " is a compiler that converts program source code into native machine code as the program is running."
$newRun1 = " is a compiler that converts program code into native machine code as the program is running."

$full = $tr.Text
$idx1 = $full.IndexOf($oldRun1)
if ($idx1 -ge 0) {
    $run1Range = $tr.Characters($idx1 + 1, $oldRun1.Length)
    $run1Range.Text = $newRun1
}

# --- Change 2 -----------------------------------------------------------
# "Execution switches to the compiled version once it becomes available,"
# is split into two runs and the trailing comma is replaced with a period:
#   "Execution switches to the compiled version once it "
#   "becomes available."
$oldTail = "becomes available,"
$newTail = "becomes available."

$full = $tr.Text
$idx2 = $full.IndexOf($oldTail)
if ($idx2 -ge 0) {
    $tailRange = $tr.Characters($idx2 + 1, $oldTail.Length)
    $tailRange.Text = $newTail
}
